# Midterm Deck edit:
#  1. Move the "Datasets" slide (last slide, slide 9) so it becomes slide 5,
#     right after "Modules Available" and before "Classification and Regression".
#  2. On the "Modules Available" slide, reorder the bullet list so that
#     "Datasets" is the first bullet (matching the new slide order) instead of
#     the last one.

$p = $ppt.ActivePresentation

# --- 1. Reorder the slides -------------------------------------------------
# The "Datasets" slide is currently the last slide (index 9); move it to
# position 5 so it directly follows "Modules Available".
$datasetsSlide = $p.Slides.Item($p.Slides.Count)
$datasetsSlide.MoveTo(5)

# --- 2. Reorder the bullets on the "Modules Available" slide ---------------
# Clear the existing text first so the new text is written fresh (avoids the
# text-diffing logic from splicing runs together when words are reused).
$modulesSlide = $p.Slides.Item(4)
$body = $modulesSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = ""
$body.Text = "Datasets`rClassification`rRegression`rTime Series`rClustering `rAnomaly Detection"
